# Generate Report for Handback
# Marks the two localization rows (zh-cn and de-de sheets) as handed back,
# records the handback target/history files + timestamp, links the new
# "Latest Target File" cell to the same source doc as column A, and widens
# a few columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1) Status column: replace "Ready for handoff" everywhere it appears
#    (Overview!E2:F3 and the Status column on the language sheets).
# ---------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange()
    foreach ($cell in $used.Cells) {
        $v = $cell.Value()
        if ($statusOld -eq $v) {
            $cell.Value = $statusNew
        }
    }
}

# ---------------------------------------------------------------------
# Helper: look up the URL behind an existing hyperlink anchored at a
# given A1 cell address on a worksheet (so the new "Latest Target File"
# link reuses exactly the same target as the "Source File Name" link).
# ---------------------------------------------------------------------
function Get-HyperlinkAddress($worksheet, $cellAddr) {
    foreach ($hl in $worksheet.Hyperlinks) {
        $rngAddr = $hl.Range().Address()
        if ($cellAddr -eq $rngAddr) {
            return $hl.Address()
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 2) Per-language-sheet handback details.
# ---------------------------------------------------------------------
$sheetInfo = @(
    @{ Name = "zh-cn"; HandbackTime = "2016-08-18 08:49:49" },
    @{ Name = "de-de"; HandbackTime = "2016-08-18 08:49:57" }
)

foreach ($info in $sheetInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Row 2 -> 7c794a4b-3e52-4c05-8f5e-59bfc11bd87f.*
    $srcName2 = $ws.Range("A2").Value()
    $url2 = Get-HyperlinkAddress $ws "`$A`$2"
    $ws.Range("I2").Value = $srcName2
    $ws.Hyperlinks.Add($ws.Range("I2"), $url2, $null, $null, $srcName2)
    $ws.Range("J2").Value = $srcName2 -replace "\.md$", (".e53680d43fa72cb1a8d826856d789cacaacbcb81." + $info.Name + ".xlf")
    $ws.Range("K2").Value = $info.HandbackTime

    # Row 3 -> e2cbc8ae-1a14-4e4b-b0f7-637666520407.*
    $srcName3 = $ws.Range("A3").Value()
    $url3 = Get-HyperlinkAddress $ws "`$A`$3"
    $ws.Range("I3").Value = $srcName3
    $ws.Hyperlinks.Add($ws.Range("I3"), $url3, $null, $null, $srcName3)
    $ws.Range("J3").Value = $srcName3 -replace "\.md$", (".c654b263278e6a1a91fbe95efde0c348a508b69f." + $info.Name + ".xlf")
    $ws.Range("K3").Value = $info.HandbackTime

    # Widen the Status / Latest Target File / Latest Handback File columns
    # now that they hold longer strings.
    $ws.Columns.Item(3).ColumnWidth = 29.14
    $ws.Columns.Item(9).ColumnWidth = 39.14
    $ws.Columns.Item(10).ColumnWidth = 39.14
}

# Overview sheet: zh-cn / de-de status columns got wider text too.
$ov = $wb.Worksheets.Item("Overview")
$ov.Columns.Item(5).ColumnWidth = 29.14
$ov.Columns.Item(6).ColumnWidth = 29.14

Write-Host "Handback report generated"
